$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 370952.88
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 370952.88
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 1112858.64
$ws.Range("N17").Value = -1113194.64
$ws.Range("H32").Value = 3140.1904
$ws.Range("I32").Value = 699.25
$ws.Range("J32").Value = 3714.5293
$ws.Range("K32").Value = 699.25
$ws.Range("L32").Value = 3714.5293
$ws.Range("M32").Value = -373.25
$ws.Range("H62").Value = 2062199.4
$ws.Range("I62").Value = 2576501.8
$ws.Range("J62").Value = 4990
$ws.Range("K62").Value = 2576501.8
$ws.Range("L62").Value = 4990
$ws.Range("M62").Value = -2575877.8
$ws.Range("N62").Value = -6238
$ws.Range("H65").Value = 2062199.4
$ws.Range("I65").Value = 2576501.8
$ws.Range("J65").Value = 4990
$ws.Range("K65").Value = 12882509
$ws.Range("L65").Value = 24950
$ws.Range("M65").Value = -12879389
$ws.Range("N65").Value = -31190
$ws.Range("H138").Value = 2767.077
$ws.Range("I138").Value = 2155.4443
$ws.Range("J138").Value = 3090.8823
$ws.Range("K138").Value = 6466.3329
$ws.Range("L138").Value = 9272.6469
$ws.Range("M138").Value = -1326.3329
$ws.Range("N138").Value = -19552.6469
$ws.Range("H141").Value = 1475
$ws.Range("I141").Value = 1475
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 4425
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 755
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1038
$ws.Range("I2").Value = 1067.6842
$ws.Range("J2").Value = 944
$ws.Range("K2").Value = 1067.6842
$ws.Range("L2").Value = 944
$ws.Range("M2").Value = -954.6841999999999
$ws.Range("H96").Value = 1000000
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 1000000
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 1000000
$ws.Range("N96").Value = -1005492
$ws.Range("H102").Value = 4368329.5
$ws.Range("I102").Value = 5051939.5
$ws.Range("J102").Value = 266670
$ws.Range("K102").Value = 5051939.5
$ws.Range("L102").Value = 266670
$ws.Range("M102").Value = -5050317.5
$ws.Range("H116").Value = 1038
$ws.Range("I116").Value = 1067.6842
$ws.Range("J116").Value = 944
$ws.Range("K116").Value = 1067.6842
$ws.Range("L116").Value = 944
$ws.Range("M116").Value = 1226.3158
$ws.Range("H132").Value = 34484684
$ws.Range("I132").Value = 40001800
$ws.Range("J132").Value = 2705
$ws.Range("K132").Value = 120005400
$ws.Range("L132").Value = 8115
$ws.Range("M132").Value = -120002870
$ws.Range("N132").Value = -13175

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1038
$ws.Range("I3").Value = 1067.6842
$ws.Range("J3").Value = 944
$ws.Range("K3").Value = 1067.6842
$ws.Range("L3").Value = 944
$ws.Range("M3").Value = -953.6841999999999
$ws.Range("H94").Value = 2688.111
$ws.Range("I94").Value = 1598.6
$ws.Range("J94").Value = 4050
$ws.Range("K94").Value = 1598.6
$ws.Range("L94").Value = 4050
$ws.Range("M94").Value = -1147.6
$ws.Range("H95").Value = 61666
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 61666
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 61666
$ws.Range("N95").Value = -67158
$ws.Range("M95").ClearContents()
$ws.Range("H105").Value = 2027.1818
$ws.Range("I105").Value = 1957
$ws.Range("J105").Value = 2150
$ws.Range("K105").Value = 1957
$ws.Range("L105").Value = 2150
$ws.Range("M105").Value = -210
$ws.Range("N105").Value = -5644
$ws.Range("H134").Value = 2548.5386
$ws.Range("I134").Value = 2455.8823
$ws.Range("J134").Value = 3178.6
$ws.Range("K134").Value = 7367.646900000001
$ws.Range("L134").Value = 9535.799999999999
$ws.Range("M134").Value = -4832.646900000001
$ws.Range("H135").Value = 58166.5
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 58166.5
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 58166.5
$ws.Range("N135").Value = -68306.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2807.116
$ws.Range("I31").Value = 2270.8948
$ws.Range("J31").Value = 3464.4194
$ws.Range("K31").Value = 2270.8948
$ws.Range("L31").Value = 3464.4194
$ws.Range("M31").Value = -1975.8948
$ws.Range("N31").Value = -4054.4194
$ws.Range("H34").Value = 2807.116
$ws.Range("I34").Value = 2270.8948
$ws.Range("J34").Value = 3464.4194
$ws.Range("K34").Value = 2270.8948
$ws.Range("L34").Value = 3464.4194
$ws.Range("M34").Value = -2068.8948
$ws.Range("N34").Value = -3868.4194
$ws.Range("H62").Value = 4728.154
$ws.Range("I62").Value = 3165.7144
$ws.Range("J62").Value = 6551
$ws.Range("K62").Value = 3165.7144
$ws.Range("L62").Value = 6551
$ws.Range("M62").Value = -2541.7144
$ws.Range("N62").Value = -7799
$ws.Range("H65").Value = 4728.154
$ws.Range("I65").Value = 3165.7144
$ws.Range("J65").Value = 6551
$ws.Range("K65").Value = 15828.572
$ws.Range("L65").Value = 32755
$ws.Range("M65").Value = -12708.572
$ws.Range("N65").Value = -38995
$ws.Range("H69").Value = 12500
$ws.Range("I69").Value = 12500
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 12500
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -11751
$ws.Range("H72").Value = 12500
$ws.Range("I72").Value = 12500
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 37500
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -33756
$ws.Range("H134").Value = 2355.4783
$ws.Range("I134").Value = 2084.4375
$ws.Range("J134").Value = 2975
$ws.Range("K134").Value = 6253.3125
$ws.Range("L134").Value = 8925
$ws.Range("M134").Value = -3718.3125
$ws.Range("N134").Value = -13995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 50114
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 50114
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 50114
$ws.Range("N63").Value = -51486
$ws.Range("H66").Value = 50114
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 50114
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 150342
$ws.Range("N66").Value = -157206
$ws.Range("H132").Value = 3716.147
$ws.Range("I132").Value = 3015.4783
$ws.Range("J132").Value = 5181.1816
$ws.Range("K132").Value = 9046.4349
$ws.Range("L132").Value = 15543.5448
$ws.Range("M132").Value = -6516.4349
$ws.Range("H136").Value = 39992.668
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 39992.668
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 119978.004
$ws.Range("N136").Value = -125078.004

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2418.625
$ws.Range("I40").Value = 2418.625
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2418.625
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -2282.625
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H99").Value = 14324.8
$ws.Range("I99").Value = 14324.8
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 14324.8
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -11329.8
$ws.Range("H136").Value = 2109.9678
$ws.Range("I136").Value = 1674.6957
$ws.Range("J136").Value = 3361.375
$ws.Range("K136").Value = 5024.0871
$ws.Range("L136").Value = 10084.125
$ws.Range("M136").Value = -2474.0871
$ws.Range("N136").Value = -15184.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 45297.75
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 45297.75
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 45297.75
$ws.Range("N46").Value = -45759.75
$ws.Range("H62").Value = 3973170
$ws.Range("I62").Value = 4766004
$ws.Range("J62").Value = 9000
$ws.Range("K62").Value = 4766004
$ws.Range("L62").Value = 9000
$ws.Range("M62").Value = -4765380
$ws.Range("H64").Value = 76591.664
$ws.Range("I64").Value = 147777
$ws.Range("J64").Value = 40999
$ws.Range("K64").Value = 147777
$ws.Range("L64").Value = 40999
$ws.Range("M64").Value = -147529
$ws.Range("N64").Value = -41495
$ws.Range("H65").Value = 3973170
$ws.Range("I65").Value = 4766004
$ws.Range("J65").Value = 9000
$ws.Range("K65").Value = 23830020
$ws.Range("L65").Value = 45000
$ws.Range("M65").Value = -23826900
$ws.Range("H67").Value = 76591.664
$ws.Range("I67").Value = 147777
$ws.Range("J67").Value = 40999
$ws.Range("K67").Value = 147777
$ws.Range("L67").Value = 40999
$ws.Range("M67").Value = -146919
$ws.Range("N67").Value = -42715
$ws.Range("H122").Value = 3848.2856
$ws.Range("I122").Value = 3864.75
$ws.Range("J122").Value = 3749.5
$ws.Range("K122").Value = 11594.25
$ws.Range("L122").Value = 11248.5
$ws.Range("M122").Value = -9144.25
$ws.Range("H126").Value = 1779.8667
$ws.Range("I126").Value = 1618
$ws.Range("J126").Value = 2225
$ws.Range("K126").Value = 4854
$ws.Range("L126").Value = 6675
$ws.Range("M126").Value = -2384
$ws.Range("N126").Value = -11615
$ws.Range("H134").Value = 45297.75
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 45297.75
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 135893.25
$ws.Range("N134").Value = -140963.25
$ws.Range("H138").Value = 45250
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 45250
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 45250
$ws.Range("N138").Value = -55530
